# EPEX Spot prices workbook update
# -----------------------------------------------------------------------
# 1) "Prix Spot" sheet : add column AE ("14-jul") with hourly prices (rows 2-25).
# 2) "Gaz" sheet        : append rows 28-29 (2025-07-12 / 2025-07-13 @ 34.8).
# 3) "CO2" sheet        : append rows 28-29 (2025-07-12 / 2025-07-13 @ 69.8).
#
# Note: literal strings that look like ISO dates ("2025-07-12") get
# auto-converted to date serials if assigned directly through .Value, the
# same way Excel itself reinterprets typed input. The source file instead
# stores those as plain text cells, so we round-trip the text through a
# helper cell's formula result and PasteSpecial(xlPasteValues) to bake in
# the literal string without triggering the date auto-detection or
# leaving any new number-format style behind.

$xlPasteValues  = -4163
$xlPasteFormats = -4122

function Set-TextValue {
    param($ws, $cellRow, $cellCol, [string]$text, $helperRow, $helperCol)

    $helper = $ws.Cells.Item($helperRow, $helperCol)
    $target = $ws.Cells.Item($cellRow, $cellCol)

    $helper.Formula = '="' + $text + '"'
    $helper.Copy()
    $target.PasteSpecial($xlPasteValues)
    $helper.ClearContents()
    $excel.CutCopyMode = $false
}

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Prix Spot" — new column AE (14-jul) ---
$wsPrix = $wb.Worksheets.Item("Prix Spot")

# Header cell AE1: "14-jul" (text, does not look like a date so .Value is safe),
# then copy the header style (bold / border / center) from AD1.
$wsPrix.Cells.Item(1, 31).Value = "14-jul"
$wsPrix.Cells.Item(1, 30).Copy()
$wsPrix.Cells.Item(1, 31).PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = $false

$aeValues = @(
    85.26000000000001,
    71.56,
    49.28,
    51.95,
    41.39,
    35.78,
    40.19,
    31.17,
    28.27,
    38.53,
    37.05,
    30.34,
    34.06,
    38.8,
    31.49,
    18.21,
    27.09,
    47.5,
    54.14,
    77.76000000000001,
    86.23,
    89.29000000000001,
    117.84,
    113.83
)

for ($i = 0; $i -lt $aeValues.Length; $i++) {
    $row = $i + 2
    $wsPrix.Cells.Item($row, 31).Value = $aeValues[$i]
}

# --- Sheet 2: "Gaz" — append rows 28 & 29 ---
$wsGaz = $wb.Worksheets.Item("Gaz")

Set-TextValue $wsGaz 28 1 "2025-07-12" 1 20
$wsGaz.Cells.Item(28, 2).Value = 34.8

Set-TextValue $wsGaz 29 1 "2025-07-13" 1 20
$wsGaz.Cells.Item(29, 2).Value = 34.8

# --- Sheet 3: "CO2" — append rows 28 & 29 ---
$wsCo2 = $wb.Worksheets.Item("CO2")

Set-TextValue $wsCo2 28 1 "2025-07-12" 1 20
$wsCo2.Cells.Item(28, 2).Value = 69.8

Set-TextValue $wsCo2 29 1 "2025-07-13" 1 20
$wsCo2.Cells.Item(29, 2).Value = 69.8
